$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("D3").Value = "데이터를 시각적으로 표현하기 위한 문법 - Grammar of graphics"
$ws.Range("E3").Value = "https://lumiamitie.github.io/data/grammar-of-graphics/"

# Row 9
$ws.Range("D9").Value = "파비클래스, PDSI를 거치지 않은 질문을 받지 말아야 할 이유 – 2"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/pabiiclass-pdsi-not-answering-2/#utm_source=rss&utm_medium=rss&utm_campaign=pabiiclass-pdsi-not-answering-2"

# Row 28
$ws.Range("D28").Value = "ROS 패키지, 스택 구조화하기 - 모바일 로봇 중심으로"
$ws.Range("E28").Value = "https://ropiens.tistory.com/146"
